# New crime data collected — weekly refresh of the 122nd Precinct CompStat
# report: bumps the report "Volume/Number" and the covered week's dates,
# and refreshes the crime-statistics table (rows 15-31) with the new
# week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: force a numeric-looking cell to hold TEXT (shared-string) data
# while keeping the original "s"-style (right aligned, General format)
# used elsewhere in the sheet for these placeholder cells ("0", "***.*").
# ---------------------------------------------------------------------
function Set-TextCell($addr, $text, $refAddr) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($refAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Header: "Volume 31   Number  25" -> "...  26"
# ---------------------------------------------------------------------
$volCell = $ws.Range("A8")
$volText = $volCell.Characters().Text
$volCell.Characters($volText.Length - 1, 2).Text = "26"

# ---------------------------------------------------------------------
# Header: "Report Covering the Week  6/17/2024  Through  6/23/2024"
#      -> "Report Covering the Week  6/24/2024  Through  6/30/2024"
# ---------------------------------------------------------------------
$weekCell = $ws.Range("C9")
$weekText = $weekCell.Characters().Text
$startIdx = $weekText.IndexOf("6/17/2024") + 1
$weekCell.Characters($startIdx, 9).Text = "6/24/2024"
$endIdx = $weekText.IndexOf("6/23/2024") + 1
$weekCell.Characters($endIdx, 9).Text = "6/30/2024"

# ---------------------------------------------------------------------
# Column H best-fit width shrinks now that its values are shorter.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 5.43

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 200
$ws.Range("L15").Value = 125

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 2
Set-TextCell "D16" "0" "D22"
Set-TextCell "E16" "***.*" "E22"
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 30
$ws.Range("K16").Value = 11.111111111111
$ws.Range("L16").Value = -3.225806451612
$ws.Range("M16").Value = -41.176470588235
$ws.Range("N16").Value = -79.020979020979

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 83
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = -11.702127659574
$ws.Range("L17").Value = 80.434782608695
$ws.Range("M17").Value = 20.289855072463
$ws.Range("N17").Value = -41.134751773049

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 44
$ws.Range("J18").Value = 49
$ws.Range("K18").Value = -10.204081632653
$ws.Range("L18").Value = 29.411764705882
$ws.Range("M18").Value = -50.561797752809
$ws.Range("N18").Value = -93.413173652694

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 36
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 259
$ws.Range("J19").Value = 238
$ws.Range("K19").Value = 8.823529411764
$ws.Range("L19").Value = 68.181818181818
$ws.Range("M19").Value = 35.602094240837
$ws.Range("N19").Value = -37.740384615384

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -30
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = -30.612244897959
$ws.Range("L20").Value = -46.875
$ws.Range("M20").Value = -22.727272727272
$ws.Range("N20").Value = -97.510980966325

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 64
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 3.225806451612
$ws.Range("I21").Value = 459
$ws.Range("J21").Value = 463
$ws.Range("K21").Value = -0.863930885529
$ws.Range("L21").Value = 37.425149700598
$ws.Range("M21").Value = 1.548672566371
$ws.Range("N21").Value = -83.290862759373

# ---------------------------------------------------------------------
# Row 23 - Petit Larceny
# ---------------------------------------------------------------------
Set-TextCell "C23" "0" "C22"
$ws.Range("L23").Value = -33.333333333333

# ---------------------------------------------------------------------
# Row 24 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -35.714285714285
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = -20.879120879120
$ws.Range("I24").Value = 514
$ws.Range("J24").Value = 562
$ws.Range("K24").Value = -8.540925266903
$ws.Range("L24").Value = 46.857142857142
$ws.Range("M24").Value = -38.954869358669

# ---------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -27.906976744186
$ws.Range("I25").Value = 246
$ws.Range("J25").Value = 246
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 215.384615384615

# ---------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 7
Set-TextCell "D26" "0" "D22"
Set-TextCell "E26" "***.*" "E22"
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 141.176470588235
$ws.Range("I26").Value = 180
$ws.Range("K26").Value = 21.621621621621
$ws.Range("L26").Value = 6.508875739644
$ws.Range("M26").Value = -39.189189189189

# ---------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 9.090909090909

# ---------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------
Set-TextCell "C28" "0" "C22"
Set-TextCell "D28" "0" "D22"
Set-TextCell "E28" "***.*" "E22"
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0

# ---------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------
Set-TextCell "G29" "0" "G22"
Set-TextCell "H29" "***.*" "H22"

# ---------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------
Set-TextCell "G30" "0" "G22"
Set-TextCell "H30" "***.*" "H22"

# ---------------------------------------------------------------------
# Row 31 - Traffic Fatalities
# ---------------------------------------------------------------------
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = -20

Write-Output "edit.ps1 applied"
